$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sample repositories (column A, rows 2-6) ---
$ws.Range("A2").Value = " GitHubSampleRepo1"
$ws.Range("A3").Value = "GitHubSampleRepo2"
$ws.Range("A4").Value = "GitHubSampleRepo3"
$ws.Range("A5").Value = " GitHubSampleRepo4"
$ws.Range("A6").Value = "GitHubSampleRepo5"

# --- Give the GitHubSampleRepo2 / 3 / 5 rows (A3, A4, A6) a distinct
#     white-fill, black Arial, left-aligned look ---
$base = $ws.Range("A3")
$base.HorizontalAlignment = -4131
$base.Interior.Color = 16777215
$base.Interior.PatternColor = 16777215
$base.Font.Color = 0
$base.Font.Name = "Arial"

$base.Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
